$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (Keyword, Correlation)
$updates = @{
    2  = @("credit", -0.2353)
    3  = @("trade", -0.2387)
    4  = @("bank", -0.1053)
    7  = @("bank", 0.4503)
    8  = @("credit", 0.0152)
    9  = @("inflation", -0.0236)
    10 = @("trade", 0.6127)
    12 = @("credit", 0.1263)
    13 = @("trade", -0.127)
    14 = @("bank", -0.1649)
    17 = @("bank", 0.0898)
    18 = @("trade", -0.1159)
    19 = @("inflation", 0.0035)
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
}
